# Issue #10 Play Playlist tidy
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Existing issue #9 (row 10) gets its Status filled in as "DONE"
$ws.Range("B10").Value = "DONE"

# New issue #10 (row 11)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "DONE"
$ws.Range("C11").Value = "UI"
# New shared strings must be registered in the same order as they appear
# in the diff: "H4 and button layout..." (F11) before "Play Playlist tidy" (D11).
$ws.Range("F11").Value = "H4 and button layout on play playlist screen"
$ws.Range("D11").Value = "Play Playlist tidy"
$ws.Range("E11").Value = "Tidy UI"

$ws.Range("D16").Select()
